$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 518 ("Poroto granado" data block),
# pushing the existing rows 518-565 down to 520-567.
$ws.Rows.Item(518).Insert()
$ws.Rows.Item(518).Insert()

# --- New row 518 ---
$ws.Cells.Item(518,1).Value = 6
$ws.Cells.Item(518,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(518,3).Value = "Metropolitana"
$ws.Cells.Item(518,4).Value = 44918
$ws.Cells.Item(518,5).Value = 13
$ws.Cells.Item(518,6).Value = 100112030
$ws.Cells.Item(518,7).Value = "Poroto granado"
$ws.Cells.Item(518,8).Value = "Sin especificar"
$ws.Cells.Item(518,9).Value = "Primera"
$ws.Cells.Item(518,10).Value = 1400
$ws.Cells.Item(518,11).Value = 30000
$ws.Cells.Item(518,12).Value = 35000
$ws.Cells.Item(518,13).Value = 32679
$ws.Cells.Item(518,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(518,15).Value = "Región Metropolitana"
$ws.Cells.Item(518,16).Value = 1307
$ws.Cells.Item(518,17).Value = 25
$ws.Cells.Item(518,18).Value = "Hortaliza"

# --- New row 519 ---
$ws.Cells.Item(519,1).Value = 6
$ws.Cells.Item(519,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(519,3).Value = "Metropolitana"
$ws.Cells.Item(519,4).Value = 44918
$ws.Cells.Item(519,5).Value = 13
$ws.Cells.Item(519,6).Value = 100112030
$ws.Cells.Item(519,7).Value = "Poroto granado"
$ws.Cells.Item(519,8).Value = "Sin especificar"
$ws.Cells.Item(519,9).Value = "Primera"
$ws.Cells.Item(519,10).Value = 470
$ws.Cells.Item(519,11).Value = 30000
$ws.Cells.Item(519,12).Value = 32000
$ws.Cells.Item(519,13).Value = 31021
$ws.Cells.Item(519,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(519,15).Value = "Región de O'Higgins"
$ws.Cells.Item(519,16).Value = 1241
$ws.Cells.Item(519,17).Value = 25
$ws.Cells.Item(519,18).Value = "Hortaliza"
